$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grades for row 7 (student #4)
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 5

# Fill in grades for row 28 (student #25)
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 5

# Update the frozen pane / active cell selection to reflect the new cursor
# position (F7) after data entry.
$ws.Range("F7").Select()
